$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -22.1755
$ws.Range("C6").Value = -12.2251
$ws.Range("A14").Value = -21.9567
$ws.Range("C18").Value = -12.2546
$ws.Range("C19").Value = -12.65110000000002
$ws.Range("A21").Value = -20.13439999999997
$ws.Range("B22").Value = 9.920299999999997
$ws.Range("A23").Value = -20.07899999999998
$ws.Range("B24").Value = 5.483400000000001
$ws.Range("A25").Value = -21.83469999999999
$ws.Range("A26").Value = -21.06669999999996
$ws.Range("B28").Value = 5.937500000000004
$ws.Range("A29").Value = -20.94949999999997
$ws.Range("B36").Value = 9.413800000000007
$ws.Range("C44").Value = -12.7855
$ws.Range("B45").Value = 4.696500000000006
$ws.Range("C47").Value = -12.24739999999999
$ws.Range("B48").Value = 5.138700000000004
$ws.Range("B49").Value = 5.663399999999995
$ws.Range("C51").Value = -10.967
$ws.Range("B52").Value = 5.411999999999995
$ws.Range("A53").Value = -21.685
$ws.Range("B53").Value = 5.744399999999996
$ws.Range("B54").Value = 4.915200000000002
$ws.Range("C55").Value = -13.48259999999999
$ws.Range("A57").Value = -21.92039999999999
$ws.Range("C57").Value = -12.69049999999999
$ws.Range("A59").Value = -22.2421
$ws.Range("C64").Value = -10.6768
$ws.Range("A69").Value = -21.69099999999997
$ws.Range("B70").Value = 5.2268
$ws.Range("A79").Value = -20.12420000000002
$ws.Range("C80").Value = -13.74250000000001
$ws.Range("A83").Value = -21.5575
$ws.Range("B86").Value = 5.429000000000002
$ws.Range("B87").Value = 5.186999999999994
$ws.Range("B89").Value = 4.475099999999996
$ws.Range("A91").Value = -20.58089999999997
$ws.Range("C92").Value = -10.4662
$ws.Range("A93").Value = -21.24300000000002
$ws.Range("C94").Value = -10.9636
$ws.Range("C96").Value = -10.1542
$ws.Range("B101").Value = 6.638299999999998
$ws.Range("C101").Value = -12.60839999999999
$ws.Range("A103").Value = -21.834
